# Applies the scheduled-runner data refresh to the per-sheet profit tables.
# Each block below corresponds to one changed row (H:N = price/profit columns).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC!row 18
$ws.Range("H18").Value = 3280  # was 2585.7144
$ws.Range("I18").Value = 3280  # was 2585.7144
$ws.Range("K18").Value = 3280  # was 2585.7144
$ws.Range("M18").Value = -2996  # was -2301.7144

# ALC!row 27
$ws.Range("H27").Value = 1072  # was 0
$ws.Range("I27").Value = 1072  # was 0
$ws.Range("K27").Value = 3216  # was 0
$ws.Range("M27").Value = -3115  # was None

# ALC!row 62
$ws.Range("H62").Value = 5250  # was 1313.75
$ws.Range("I62").Value = 0  # was 1313.75
$ws.Range("J62").Value = 5250  # was 0
$ws.Range("K62").Value = 0  # was 1313.75
$ws.Range("L62").Value = 5250  # was 0
$ws.Range("M62").ClearContents()  # was -689.75
$ws.Range("N62").Value = -6498  # was None

# ALC!row 65
$ws.Range("H65").Value = 5250  # was 1313.75
$ws.Range("I65").Value = 0  # was 1313.75
$ws.Range("J65").Value = 5250  # was 0
$ws.Range("K65").Value = 0  # was 6568.75
$ws.Range("L65").Value = 26250  # was 0
$ws.Range("M65").ClearContents()  # was -3448.75
$ws.Range("N65").Value = -32490  # was None

# ALC!row 81
$ws.Range("H81").Value = 60000  # was 0
$ws.Range("J81").Value = 60000  # was 0
$ws.Range("L81").Value = 60000  # was 0
$ws.Range("N81").Value = -61996  # was None

# ALC!row 84
$ws.Range("H84").Value = 60000  # was 0
$ws.Range("J84").Value = 60000  # was 0
$ws.Range("L84").Value = 180000  # was 0
$ws.Range("N84").Value = -189984  # was None

# ALC!row 92
$ws.Range("H92").Value = 309  # was 376.55554
$ws.Range("I92").Value = 207.66667  # was 236.125
$ws.Range("J92").Value = 765  # was 1500
$ws.Range("K92").Value = 207.66667  # was 236.125
$ws.Range("L92").Value = 765  # was 1500
$ws.Range("M92").Value = 1040.33333  # was 1011.875
$ws.Range("N92").Value = -3261  # was -3996

# ALC!row 98
$ws.Range("H98").Value = 758.5714  # was 786.3077
$ws.Range("J98").Value = 499  # was 549.5
$ws.Range("L98").Value = 499  # was 549.5
$ws.Range("N98").Value = -3495  # was -3545.5

# ALC!row 99
$ws.Range("H99").Value = 850  # was 916.6667
$ws.Range("I99").Value = 300  # was 375
$ws.Range("J99").Value = 1125  # was 2000
$ws.Range("K99").Value = 900  # was 1125
$ws.Range("L99").Value = 3375  # was 6000
$ws.Range("M99").Value = 598  # was 373
$ws.Range("N99").Value = -6371  # was -8996

# ALC!row 101
$ws.Range("H101").Value = 8334560.5  # was 10003014
$ws.Range("I101").Value = 20000398  # was 16669815
$ws.Range("J101").Value = 1819.8572  # was 2812.25
$ws.Range("K101").Value = 60001194  # was 50009445
$ws.Range("L101").Value = 5459.571599999999  # was 8436.75
$ws.Range("M101").Value = -59999572  # was -50007823
$ws.Range("N101").Value = -8703.571599999999  # was -11680.75

# ALC!row 116
$ws.Range("H116").Value = 4057.3333  # was 4134.2856
$ws.Range("I116").Value = 3240  # was 3326.6667
$ws.Range("K116").Value = 3240  # was 3326.6667
$ws.Range("M116").Value = 202  # was 115.3332999999998

# ALC!row 118
$ws.Range("H118").Value = 1145.625  # was 1062.7778
$ws.Range("I118").Value = 833.4  # was 761.1667
$ws.Range("K118").Value = 2500.2  # was 2283.5001
$ws.Range("M118").Value = -843.1999999999998  # was -626.5001000000002

# ALC!row 122
$ws.Range("H122").Value = 758.5714  # was 786.3077
$ws.Range("J122").Value = 499  # was 549.5
$ws.Range("L122").Value = 1497  # was 1648.5
$ws.Range("N122").Value = -6397  # was -6548.5

# ALC!row 135
$ws.Range("H135").Value = 1989.6  # was 1221.9445
$ws.Range("I135").Value = 1877.4445  # was 1291.0714
$ws.Range("J135").Value = 2999  # was 980
$ws.Range("K135").Value = 16897.0005  # was 11619.6426
$ws.Range("L135").Value = 26991  # was 8820
$ws.Range("M135").Value = -14362.0005  # was -9084.642600000001
$ws.Range("N135").Value = -32061  # was -13890

# ALC!row 137
$ws.Range("H137").Value = 1909.5714  # was 1939.8
$ws.Range("I137").Value = 1928.4  # was 1983
$ws.Range("J137").Value = 1862.5  # was 1875
$ws.Range("K137").Value = 5785.200000000001  # was 5949
$ws.Range("L137").Value = 5587.5  # was 5625
$ws.Range("M137").Value = -3235.200000000001  # was -3399
$ws.Range("N137").Value = -10687.5  # was -10725

# ALC!row 141
$ws.Range("H141").Value = 2523.1538  # was 2695.9167
$ws.Range("I141").Value = 2402.111  # was 2646.125
$ws.Range("K141").Value = 7206.333  # was 7938.375
$ws.Range("M141").Value = -2026.333  # was -2758.375

$ws = $wb.Worksheets.Item("ARM")
# ARM!row 45
$ws.Range("H45").Value = 2997.2693  # was 2752.96
$ws.Range("I45").Value = 2112.7778  # was 2043.421
$ws.Range("J45").Value = 4987.375  # was 4999.8335
$ws.Range("K45").Value = 2112.7778  # was 2043.421
$ws.Range("L45").Value = 4987.375  # was 4999.8335
$ws.Range("M45").Value = -1735.7778  # was -1666.421
$ws.Range("N45").Value = -5741.375  # was -5753.8335

# ARM!row 46
$ws.Range("H46").Value = 9232.333000000001  # was 9223
$ws.Range("J46").Value = 9200  # was 9197.5
$ws.Range("L46").Value = 9200  # was 9197.5
$ws.Range("N46").Value = -9838  # was -9835.5

# ARM!row 61
$ws.Range("H61").Value = 6487.5  # was 7061.769
$ws.Range("I61").Value = 2972.7144  # was 3134.8333
$ws.Range("J61").Value = 9221.223  # was 10427.714
$ws.Range("K61").Value = 2972.7144  # was 3134.8333
$ws.Range("L61").Value = 9221.223  # was 10427.714
$ws.Range("M61").Value = -2760.7144  # was -2922.8333
$ws.Range("N61").Value = -9645.223  # was -10851.714

# ARM!row 102
$ws.Range("H102").Value = 1845.6  # was 1877.6
$ws.Range("I102").Value = 1576  # was 1629.3334
$ws.Range("K102").Value = 1576  # was 1629.3334
$ws.Range("M102").Value = 46  # was -7.333399999999983

# ARM!row 122
$ws.Range("H122").Value = 13789.5  # was 18017.63
$ws.Range("I122").Value = 18745.285  # was 20934.812
$ws.Range("J122").Value = 2226  # was 2459.3333
$ws.Range("K122").Value = 56235.855  # was 62804.436
$ws.Range("L122").Value = 6678  # was 7377.999899999999
$ws.Range("M122").Value = -53785.855  # was -60354.436
$ws.Range("N122").Value = -11578  # was -12277.9999

# ARM!row 136
$ws.Range("H136").Value = 6487.5  # was 7061.769
$ws.Range("I136").Value = 2972.7144  # was 3134.8333
$ws.Range("J136").Value = 9221.223  # was 10427.714
$ws.Range("K136").Value = 8918.143199999999  # was 9404.499899999999
$ws.Range("L136").Value = 27663.669  # was 31283.142
$ws.Range("M136").Value = -6368.143199999999  # was -6854.499899999999
$ws.Range("N136").Value = -32763.669  # was -36383.142

$ws = $wb.Worksheets.Item("BSM")
# BSM!row 20
$ws.Range("H20").Value = 3620.4  # was 3125.3333
$ws.Range("J20").Value = 2597  # was 1948
$ws.Range("L20").Value = 2597  # was 1948
$ws.Range("N20").Value = -3091  # was -2442

# BSM!row 63
$ws.Range("H63").Value = 35000  # was 0
$ws.Range("J63").Value = 35000  # was 0
$ws.Range("L63").Value = 35000  # was 0
$ws.Range("N63").Value = -36372  # was None

# BSM!row 66
$ws.Range("H66").Value = 35000  # was 0
$ws.Range("J66").Value = 35000  # was 0
$ws.Range("L66").Value = 105000  # was 0
$ws.Range("N66").Value = -111864  # was None

$ws = $wb.Worksheets.Item("CRP")
# CRP!row 31
$ws.Range("H31").Value = 963.6923  # was 1230
$ws.Range("I31").Value = 921.8182  # was 960
$ws.Range("J31").Value = 1194  # was 1500
$ws.Range("K31").Value = 921.8182  # was 960
$ws.Range("L31").Value = 1194  # was 1500
$ws.Range("M31").Value = -626.8182  # was -665
$ws.Range("N31").Value = -1784  # was -2090

# CRP!row 34
$ws.Range("H34").Value = 963.6923  # was 1230
$ws.Range("I34").Value = 921.8182  # was 960
$ws.Range("J34").Value = 1194  # was 1500
$ws.Range("K34").Value = 921.8182  # was 960
$ws.Range("L34").Value = 1194  # was 1500
$ws.Range("M34").Value = -719.8182  # was -758
$ws.Range("N34").Value = -1598  # was -1904

# CRP!row 74
$ws.Range("H74").Value = 75157  # was 100314
$ws.Range("J74").Value = 75157  # was 100314
$ws.Range("L74").Value = 75157  # was 100314
$ws.Range("N74").Value = -76905  # was -102062

# CRP!row 77
$ws.Range("H77").Value = 75157  # was 100314
$ws.Range("J77").Value = 75157  # was 100314
$ws.Range("L77").Value = 225471  # was 300942
$ws.Range("N77").Value = -234207  # was -309678

$ws = $wb.Worksheets.Item("CUL")
# CUL!row 4
$ws.Range("H4").Value = 9239597  # was 8007676.5
$ws.Range("I4").Value = 10009147  # was 8579296
$ws.Range("K4").Value = 30027441  # was 25737888
$ws.Range("M4").Value = -30027329  # was -25737776

# CUL!row 23
$ws.Range("H23").Value = 718.5  # was 747.86664
$ws.Range("I23").Value = 597.5  # was 583.5714
$ws.Range("J23").Value = 839.5  # was 891.625
$ws.Range("K23").Value = 1792.5  # was 1750.7142
$ws.Range("L23").Value = 2518.5  # was 2674.875
$ws.Range("M23").Value = -1557.5  # was -1515.7142
$ws.Range("N23").Value = -2988.5  # was -3144.875

# CUL!row 68
$ws.Range("H68").Value = 2704.761  # was 2715.0222
$ws.Range("J68").Value = 2852.524  # was 2867.3901
$ws.Range("L68").Value = 8557.572  # was 8602.1703
$ws.Range("N68").Value = -10179.572  # was -10224.1703

# CUL!row 71
$ws.Range("H71").Value = 2704.761  # was 2715.0222
$ws.Range("J71").Value = 2852.524  # was 2867.3901
$ws.Range("L71").Value = 25672.716  # was 25806.5109
$ws.Range("N71").Value = -33784.716  # was -33918.5109

# CUL!row 92
$ws.Range("H92").Value = 2250  # was 1500
$ws.Range("J92").Value = 3000  # was 0
$ws.Range("L92").Value = 9000  # was 0
$ws.Range("N92").Value = -11496  # was None

# CUL!row 121
$ws.Range("H121").Value = 9914.5  # was 9921.0625
$ws.Range("J121").Value = 4608.1  # was 4618.6
$ws.Range("L121").Value = 13824.3  # was 13855.8
$ws.Range("N121").Value = -16444.3  # was -16475.8

$ws = $wb.Worksheets.Item("GSM")
# GSM!row 36
$ws.Range("H36").Value = 12880285  # was 22519750
$ws.Range("I36").Value = 27998.75  # was 29000
$ws.Range("K36").Value = 27998.75  # was 29000
$ws.Range("M36").Value = -27513.75  # was -28515

# GSM!row 102
$ws.Range("H102").Value = 2198  # was 2067.2
$ws.Range("I102").Value = 2198  # was 2993
$ws.Range("J102").Value = 0  # was 1450
$ws.Range("K102").Value = 2198  # was 2993
$ws.Range("L102").Value = 0  # was 1450
$ws.Range("M102").Value = -576  # was -1371
$ws.Range("N102").ClearContents()  # was -4694

$ws = $wb.Worksheets.Item("LTW")
# LTW!row 132
$ws.Range("H132").Value = 2840.1177  # was 2949.0715
$ws.Range("I132").Value = 2937.2307  # was 3053.5454
$ws.Range("J132").Value = 2524.5  # was 2566
$ws.Range("K132").Value = 8811.6921  # was 9160.636200000001
$ws.Range("L132").Value = 7573.5  # was 7698
$ws.Range("M132").Value = -6281.6921  # was -6630.636200000001
$ws.Range("N132").Value = -12633.5  # was -12758

$ws = $wb.Worksheets.Item("WVR")
# WVR!row 122
$ws.Range("H122").Value = 4292.8887  # was 3998.6
$ws.Range("I122").Value = 2358.3333  # was 2106.25
$ws.Range("K122").Value = 7074.999899999999  # was 6318.75
$ws.Range("M122").Value = -4624.999899999999  # was -3868.75

# WVR!row 132
$ws.Range("H132").Value = 1367.5333  # was 1718.25
$ws.Range("I132").Value = 1411.6428  # was 1806.2727
$ws.Range("K132").Value = 4234.928400000001  # was 5418.8181
$ws.Range("M132").Value = -1704.928400000001  # was -2888.8181

# WVR!row 136
$ws.Range("H136").Value = 2128.0715  # was 2318.7307
$ws.Range("I136").Value = 2140.037  # was 2339.28
$ws.Range("K136").Value = 6420.110999999999  # was 7017.84
$ws.Range("M136").Value = -3870.110999999999  # was -4467.84
